$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Delete() | Out-Null
$ws.Rows.Item(5).Select() | Out-Null
